# Applies the per-cell profit/price recalculation updates produced by the
# scheduled market-data refresh run across all eight job sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 600919.4
$ws.Range("I33").Value = 772
$ws.Range("J33").Value = 1561155.2
$ws.Range("K33").Value = 772
$ws.Range("L33").Value = 1561155.2
$ws.Range("M33").Value = -543
$ws.Range("N33").Value = -1561613.2
$ws.Range("H138").Value = 2580.1648
$ws.Range("I138").Value = 1475.7826
$ws.Range("J138").Value = 2989.8547
$ws.Range("K138").Value = 4427.3478
$ws.Range("L138").Value = 8969.5641
$ws.Range("M138").Value = 712.6522000000004
$ws.Range("N138").Value = -19249.5641

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 32270.125
$ws.Range("I2").Value = 1056.7241
$ws.Range("K2").Value = 1056.7241
$ws.Range("M2").Value = -943.7240999999999
$ws.Range("H32").Value = 4325.24
$ws.Range("I32").Value = 3646.848
$ws.Range("J32").Value = 12126.75
$ws.Range("K32").Value = 3646.848
$ws.Range("L32").Value = 12126.75
$ws.Range("M32").Value = -3359.848
$ws.Range("N32").Value = -12700.75
$ws.Range("H74").Value = 844.3871
$ws.Range("I74").Value = 771.5417
$ws.Range("J74").Value = 1094.1428
$ws.Range("K74").Value = 771.5417
$ws.Range("L74").Value = 1094.1428
$ws.Range("M74").Value = 102.4583
$ws.Range("N74").Value = -2842.1428
$ws.Range("H77").Value = 844.3871
$ws.Range("I77").Value = 771.5417
$ws.Range("J77").Value = 1094.1428
$ws.Range("K77").Value = 3857.7085
$ws.Range("L77").Value = 5470.714
$ws.Range("M77").Value = 510.2915000000003
$ws.Range("N77").Value = -14206.714
$ws.Range("H110").Value = 71578980
$ws.Range("I110").Value = 71578980
$ws.Range("K110").Value = 71578980
$ws.Range("M110").Value = -71576935
$ws.Range("H116").Value = 32270.125
$ws.Range("I116").Value = 1056.7241
$ws.Range("K116").Value = 1056.7241
$ws.Range("M116").Value = 1237.2759
$ws.Range("H122").Value = 1966.2174
$ws.Range("I122").Value = 1712.15
$ws.Range("J122").Value = 3660
$ws.Range("K122").Value = 5136.450000000001
$ws.Range("L122").Value = 10980
$ws.Range("M122").Value = -2686.450000000001
$ws.Range("N122").Value = -15880
$ws.Range("H132").Value = 26377.912
$ws.Range("I132").Value = 30662.264
$ws.Range("J132").Value = 6027.25
$ws.Range("K132").Value = 91986.792
$ws.Range("L132").Value = 18081.75
$ws.Range("M132").Value = -89456.792
$ws.Range("N132").Value = -23141.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 32270.125
$ws.Range("I3").Value = 1056.7241
$ws.Range("K3").Value = 1056.7241
$ws.Range("M3").Value = -942.7240999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 34511.34
$ws.Range("I31").Value = 667.26666
$ws.Range("J31").Value = 52016.9
$ws.Range("K31").Value = 667.26666
$ws.Range("L31").Value = 52016.9
$ws.Range("M31").Value = -372.26666
$ws.Range("N31").Value = -52606.9
$ws.Range("H34").Value = 34511.34
$ws.Range("I34").Value = 667.26666
$ws.Range("J34").Value = 52016.9
$ws.Range("K34").Value = 667.26666
$ws.Range("L34").Value = 52016.9
$ws.Range("M34").Value = -465.26666
$ws.Range("N34").Value = -52420.9
$ws.Range("H58").Value = 1410.8511
$ws.Range("I58").Value = 1243.6296
$ws.Range("J58").Value = 1636.6
$ws.Range("K58").Value = 1243.6296
$ws.Range("L58").Value = 1636.6
$ws.Range("M58").Value = -1040.6296
$ws.Range("N58").Value = -2042.6
$ws.Range("H99").Value = 2562.7334
$ws.Range("I99").Value = 2123
$ws.Range("J99").Value = 2722.6365
$ws.Range("K99").Value = 2123
$ws.Range("L99").Value = 2722.6365
$ws.Range("M99").Value = -625
$ws.Range("N99").Value = -5718.636500000001
$ws.Range("H126").Value = 2562.7334
$ws.Range("I126").Value = 2123
$ws.Range("J126").Value = 2722.6365
$ws.Range("K126").Value = 6369
$ws.Range("L126").Value = 8167.9095
$ws.Range("M126").Value = -3899
$ws.Range("N126").Value = -13107.9095
$ws.Range("H132").Value = 53574110
$ws.Range("I132").Value = 55558624
$ws.Range("J132").Value = 50001990
$ws.Range("K132").Value = 166675872
$ws.Range("L132").Value = 150005970
$ws.Range("M132").Value = -166673342
$ws.Range("N132").Value = -150011030
$ws.Range("H134").Value = 1409.0883
$ws.Range("I134").Value = 797.5599999999999
$ws.Range("J134").Value = 3107.7778
$ws.Range("K134").Value = 2392.68
$ws.Range("L134").Value = 9323.3334
$ws.Range("M134").Value = 142.3200000000002
$ws.Range("N134").Value = -14393.3334
$ws.Range("H136").Value = 1410.8511
$ws.Range("I136").Value = 1243.6296
$ws.Range("J136").Value = 1636.6
$ws.Range("K136").Value = 3730.8888
$ws.Range("L136").Value = 4909.799999999999
$ws.Range("M136").Value = -1180.8888
$ws.Range("N136").Value = -10009.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 781.3077
$ws.Range("I2").Value = 11.666667
$ws.Range("J2").Value = 1012.2
$ws.Range("K2").Value = 70.00000199999999
$ws.Range("L2").Value = 6073.200000000001
$ws.Range("M2").Value = 42.99999800000001
$ws.Range("N2").Value = -6299.200000000001
$ws.Range("H58").Value = 2100
$ws.Range("J58").Value = 1400
$ws.Range("L58").Value = 4200
$ws.Range("N58").Value = -4456

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1988.8422
$ws.Range("I122").Value = 1543.3334
$ws.Range("K122").Value = 4630.0002
$ws.Range("M122").Value = -2180.0002
$ws.Range("H126").Value = 3430.7856
$ws.Range("I126").Value = 3642.75
$ws.Range("J126").Value = 3148.1667
$ws.Range("K126").Value = 10928.25
$ws.Range("L126").Value = 9444.500100000001
$ws.Range("M126").Value = -8458.25
$ws.Range("N126").Value = -14384.5001
$ws.Range("H132").Value = 2850.7666
$ws.Range("I132").Value = 2023.8
$ws.Range("J132").Value = 4504.7
$ws.Range("K132").Value = 6071.4
$ws.Range("L132").Value = 13514.1
$ws.Range("M132").Value = -3541.4
$ws.Range("N132").Value = -18574.1

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 3638.5557
$ws.Range("I93").Value = 4240.75
$ws.Range("J93").Value = 3156.8
$ws.Range("K93").Value = 4240.75
$ws.Range("L93").Value = 3156.8
$ws.Range("M93").Value = -2992.75
$ws.Range("N93").Value = -5652.8
$ws.Range("H122").Value = 3003.6191
$ws.Range("I122").Value = 2786.125
$ws.Range("J122").Value = 3699.6
$ws.Range("K122").Value = 8358.375
$ws.Range("L122").Value = 11098.8
$ws.Range("M122").Value = -5908.375
$ws.Range("N122").Value = -15998.8
$ws.Range("H132").Value = 7016.6924
$ws.Range("J132").Value = 7277.6
$ws.Range("L132").Value = 21832.8
$ws.Range("N132").Value = -26892.8
$ws.Range("H133").Value = 35730.4
$ws.Range("J133").Value = 35730.4
$ws.Range("L133").Value = 35730.4
$ws.Range("N133").Value = -40790.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 167356.33
$ws.Range("I107").Value = 848
$ws.Range("J107").Value = 333864.66
$ws.Range("K107").Value = 2544
$ws.Range("L107").Value = 1001593.98
$ws.Range("M107").Value = -624
$ws.Range("N107").Value = -1005433.98
$ws.Range("H113").Value = 751.25
$ws.Range("I113").Value = 303.25
$ws.Range("J113").Value = 975.25
$ws.Range("K113").Value = 909.75
$ws.Range("L113").Value = 2925.75
$ws.Range("M113").Value = 1260.25
$ws.Range("N113").Value = -7265.75
$ws.Range("H122").Value = 2238.2856
$ws.Range("I122").Value = 1745.421
$ws.Range("J122").Value = 3278.7778
$ws.Range("K122").Value = 5236.263
$ws.Range("L122").Value = 9836.3334
$ws.Range("M122").Value = -2786.263
$ws.Range("N122").Value = -14736.3334
